# Trading update: 2026-02-17 08:08:57
# Appends a new closed/open trade row (row 13) to both the "All Trades"
# and "MarketMaking" worksheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("All Trades", "MarketMaking")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $row = 13

    $ws.Cells.Item($row, 1).Value = 12

    # Dates/times are stored as plain text in this workbook, so force a
    # text number format first to stop Excel auto-converting the literal
    # "2026-02-17" / "08:08:51" strings into date/time serial numbers.
    $ws.Cells.Item($row, 2).NumberFormat = "@"
    $ws.Cells.Item($row, 2).Value = "2026-02-17"
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = "08:08:51"

    $ws.Cells.Item($row, 4).Value = "MarketMaking"
    $ws.Cells.Item($row, 5).Value = "DOWN"
    $ws.Cells.Item($row, 6).Value = 0.53

    # Exit Price is blank (empty string) while the trade is still OPEN.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = ""

    $ws.Cells.Item($row, 8).Value = "OPEN"
    $ws.Cells.Item($row, 9).Value = 0
    $ws.Cells.Item($row, 10).Value = 0
    $ws.Cells.Item($row, 11).Value = 99.8206335792346
    $ws.Cells.Item($row, 12).Value = 0
    $ws.Cells.Item($row, 13).Value = 0
    $ws.Cells.Item($row, 14).Value = 0.6
    $ws.Cells.Item($row, 15).Value = "Normal spread capture: 19600 bps"

    # Exit Reason is blank (empty string) while the trade is still OPEN.
    $ws.Cells.Item($row, 16).NumberFormat = "@"
    $ws.Cells.Item($row, 16).Value = ""

    $ws.Cells.Item($row, 17).Value = 0
}
